# Update Jane Smith's Salary from 0 to 100 on the Employees sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 100
